$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "114091a778-1a50787db8"
$ws.Range("B2").Value = "114091a778-9066f1e71b"

$ws.Range("A3").Value = "196d11e6b8-341ce94b5d"
$ws.Range("B3").Value = "196d11e6b8-a6b5cb57a3"

$ws.Range("A4").Value = "19fb0e96dd-63b09e56e9"
$ws.Range("B4").Value = "19fb0e96dd-63b09e56e9"

$ws.Range("A5").Value = "48aba86cb8-702980d58a"
$ws.Range("B5").Value = "48aba86cb8-702980d58a"

$ws.Range("A6").Value = "7d066011f6-6a35e27a02"
$ws.Range("B6").Value = "7d066011f6-6a35e27a02"

$ws.Range("A7").Value = "8a8394d9eb-a4fafb4123"
$ws.Range("B7").Value = "8a8394d9eb-a4fafb4123"

$ws.Range("A8").Value = "d5722a624b-ad620af736"
$ws.Range("B8").Value = "d5722a624b-ad620af736"

$ws.Range("A9").Value = "d8499f5e39-aa576e641b"
$ws.Range("B9").Value = "d8499f5e39-aa576e641b"

$ws.Range("A10").Value = "e2e1425a6e-665ca36a3e"
$ws.Range("B10").Value = "e2e1425a6e-665ca36a3e"
